# Delete the row for "PAUTA ASTUDILLO JULIO HERNAN" (row 39) on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. Excel automatically shifts
# the remaining rows up and keeps the trailing totals row anchored at the
# (now one-row-higher) bottom of the range.

$wb = $excel.ActiveWorkbook

$sheetNames = @("VENTAS POR GRUPO", "VENTA MENSUAL")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(39).Delete()
}

# The "VENTAS POR GRUPO" sheet has a trailing "N de 54" style summary row
# (counts of non-zero entries out of the 54 clients). Since one client row
# was removed, the denominator must be updated from 54 to 53. These are
# static text values (not formulas), so they are fixed up explicitly here
# now that the row has shifted up to row 55.
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(55, $col)
    $cell.Value = $cell.Text -replace "de 54", "de 53"
}
